$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 45912
$ws.Range("B2").Value = 5913.14141862252
$ws.Range("C2").Value = 4464.08865436985
$ws.Range("D2").Value = 5112
$ws.Range("E2").Value = 6520.665898
$ws.Range("F2").Value = -1.68278609386152

# Row 3
$ws.Range("A3").Value = 45913
$ws.Range("B3").Value = 1177.4817126383
$ws.Range("C3").Value = 2210.1312608157
$ws.Range("D3").Value = 2952
$ws.Range("E3").Value = 2248.969119
$ws.Range("F3").Value = 13.7341111323916

# Row 4
$ws.Range("A4").Value = 45914
$ws.Range("B4").Value = 1060.75017668403
$ws.Range("C4").Value = 2165.93148897365
$ws.Range("D4").Value = 2952
$ws.Range("E4").Value = 2127.948681
$ws.Range("F4").Value = 11.7137497204009

# Row 5
$ws.Range("A5").Value = 45915
$ws.Range("B5").Value = 5883.51545958832
$ws.Range("C5").Value = 4942.23846518091
$ws.Range("D5").Value = 2952
$ws.Range("E5").Value = 6537.637225
$ws.Range("F5").Value = 110.181676274691

# Row 6
$ws.Range("A6").Value = 45916
$ws.Range("B6").Value = 5880.24471902051
$ws.Range("C6").Value = 5132.52330834507
$ws.Range("D6").Value = 2952
$ws.Range("E6").Value = 6534.077719
$ws.Range("F6").Value = 118.098179513524

# Row 7
$ws.Range("A7").Value = 45917
$ws.Range("B7").Value = 6061.80298841384
$ws.Range("C7").Value = 5244.51158207834
$ws.Range("D7").Value = 2952
$ws.Range("E7").Value = 6722.442273
$ws.Range("F7").Value = 123.047952777688

# Row 8
$ws.Range("A8").Value = 45918
$ws.Range("B8").Value = 6061.80298841384
$ws.Range("C8").Value = 5235.62105586547
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 6722.442273
$ws.Range("F8").Value = 122.677514185485

# Row 9
$ws.Range("A9").Value = 45919
$ws.Range("B9").Value = 6061.80298841384
$ws.Range("C9").Value = 4518.21451058875
$ws.Range("D9").Value = 2952
$ws.Range("E9").Value = 6722.442273
$ws.Range("F9").Value = 92.7855747989547

# Row 10
$ws.Range("A10").Value = 45920
$ws.Range("B10").Value = 1187.6944176276
$ws.Range("C10").Value = 2255.44951038635
$ws.Range("D10").Value = 2952
$ws.Range("E10").Value = 2248.907555
$ws.Range("F10").Value = 15.1942769899479

# Row 11
$ws.Range("A11").Value = 45921
$ws.Range("B11").Value = 1056.51488323291
$ws.Range("C11").Value = 2230.35898897121
$ws.Range("D11").Value = 2952
$ws.Range("E11").Value = 2106.590982
$ws.Range("F11").Value = 13.6847953224291

# Row 12
$ws.Range("A12").Value = 45922
$ws.Range("B12").Value = 5979.99309220113
$ws.Range("C12").Value = 5075.35714801183
$ws.Range("D12").Value = 2952
$ws.Range("E12").Value = 6713.303219
$ws.Range("F12").Value = 119.027803117112

# Row 13
$ws.Range("A13").Value = 45923
$ws.Range("B13").Value = 5979.99309220113
$ws.Range("C13").Value = 5422.37678388834
$ws.Range("D13").Value = 2952
$ws.Range("E13").Value = 6713.303219
$ws.Range("F13").Value = 133.486954611967

# Row 14
$ws.Range("A14").Value = 45924
$ws.Range("B14").Value = 5979.99309220113
$ws.Range("C14").Value = 6002.14974943777
$ws.Range("D14").Value = 2952
$ws.Range("E14").Value = 6713.303219
$ws.Range("F14").Value = 157.64416150986

# Row 15
$ws.Range("A15").Value = 45925
$ws.Range("B15").Value = 5979.99309220113
$ws.Range("C15").Value = 5947.77417027515
$ws.Range("D15").Value = 2952
$ws.Range("E15").Value = 6713.303219
$ws.Range("F15").Value = 155.378512378084

